# Update non NMA parameter estimates
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("state_utility")
$ws2 = $wb.Worksheets.Item("ae_disutility")

# ---------------------------------------------------------------------------
# 1. New shared strings, created in the same order as the original edit so
#    the shared-string table indices line up with the source workbook.
# ---------------------------------------------------------------------------
$ws1.Range("D2").Value = "nafees2017health"

$ws2.Range("F1").Value = "notes"
$ws2.Range("F5").Value = "Assumed equal to rash disutility"
$ws2.Range("F6").Value = "Assumed equal to hair loss"
$ws2.Range("F8").Value = "Assumed to dyspnea"
$ws2.Range("F2").Value = "Assumed equal to fatigue disutility"
$ws2.Range("E8").Value = "doyle2008health"

# ---------------------------------------------------------------------------
# 2. state_utility sheet - updated parameter estimates
# ---------------------------------------------------------------------------
$ws1.Range("B2").Value = 0.754
$ws1.Range("C2").Value = 0
$ws1.Range("C2").Interior.Color = 65535

# ---------------------------------------------------------------------------
# 3. ae_disutility sheet - fill in mean/se/ref/notes columns
# ---------------------------------------------------------------------------
$ws2.Range("C2").Value = 0.07346
$ws2.Range("D2").Value = 0.01849
$ws2.Range("E2").Value = "nafees2008health"

$ws2.Range("C3").Value = 0.07346
$ws2.Range("D3").Value = 0.01849
$ws2.Range("E3").Value = "nafees2008health"
$ws2.Range("F3").Value = "Assumed equal to fatigue disutility"

$ws2.Range("E4").Value = "nafees2008health"

$ws2.Range("C5").Value = 0.03248
$ws2.Range("D5").Value = 0.01171
$ws2.Range("E5").Value = "nafees2008health"

$ws2.Range("C6").Value = 0.04495
$ws2.Range("D6").Value = 0.01482
$ws2.Range("E6").Value = "nafees2008health"

$ws2.Range("C7").Value = 0.03248
$ws2.Range("D7").Value = 0.01171
$ws2.Range("E7").Value = "nafees2008health"
$ws2.Range("F7").Value = "Assumed equal to rash disutility"

$ws2.Range("C8").Value = 0.05
$ws2.Range("D8").Value = 0.012

$ws2.Range("C9").Value = 0.03248
$ws2.Range("D9").Value = 0.01171
$ws2.Range("E9").Value = "nafees2008health"
$ws2.Range("F9").Value = "Assumed equal to rash disutility"

$ws2.Range("E10").Value = "nafees2008health"

$ws2.Range("C11").Value = 0.03248
$ws2.Range("D11").Value = 0.01171
$ws2.Range("E11").Value = "nafees2008health"
$ws2.Range("F11").Value = "Assumed equal to rash disutility"

# column E (ref) grew new entries - width recalculated by the authoring tool
$ws2.Columns.Item(5).ColumnWidth = 16

# G8 - touched (empty) cell that widens the sheet's used range to G11
$ws2.Range("G8").Interior.Pattern = -4142

# ---------------------------------------------------------------------------
# 4. Selections / active sheet - ae_disutility becomes the active tab
# ---------------------------------------------------------------------------
$ws1.Range("D2").Select() | Out-Null
$ws2.Activate() | Out-Null
$ws2.Range("E8").Select() | Out-Null

Write-Host "edit complete"
